$d = $word.ActiveDocument

# Collapse a range to the very end of the document body (after the last
# paragraph, "Archivo de ejemplo para curso de git.").
$e = $d.Content.End
$r = $d.Range($e, $e)

# Build a WordprocessingML package fragment containing the three new
# paragraphs to append:
#   1) an empty paragraph
#   2) "Segunda linea en word" (with the spell-check proofing marks Word
#      itself adds around each word it does not recognize)
#   3) a trailing empty paragraph
$q = [char]34
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$pkgNs = "http://schemas.microsoft.com/office/2006/xmlPackage"

$newParas = ""
$newParas += "<w:p/>"
$newParas += "<w:p>"
$newParas += "<w:proofErr w:type=" + $q + "spellStart" + $q + "/>"
$newParas += "<w:r><w:t>Segunda</w:t></w:r>"
$newParas += "<w:proofErr w:type=" + $q + "spellEnd" + $q + "/>"
$newParas += "<w:r><w:t xml:space=" + $q + "preserve" + $q + "> </w:t></w:r>"
$newParas += "<w:proofErr w:type=" + $q + "spellStart" + $q + "/>"
$newParas += "<w:r><w:t>linea</w:t></w:r>"
$newParas += "<w:proofErr w:type=" + $q + "spellEnd" + $q + "/>"
$newParas += "<w:r><w:t xml:space=" + $q + "preserve" + $q + "> en </w:t></w:r>"
$newParas += "<w:proofErr w:type=" + $q + "spellStart" + $q + "/>"
$newParas += "<w:r><w:t>word</w:t></w:r>"
$newParas += "<w:proofErr w:type=" + $q + "spellEnd" + $q + "/>"
$newParas += "</w:p>"
$newParas += "<w:p/>"

$xmlFragment = "<pkg:package xmlns:pkg=" + $q + $pkgNs + $q + ">" +
    "<pkg:part pkg:name=" + $q + "/word/document.xml" + $q + " pkg:contentType=" + $q + "application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" + $q + ">" +
    "<pkg:xmlData>" +
    "<w:document xmlns:w=" + $q + $wNs + $q + ">" +
    "<w:body>" + $newParas + "</w:body>" +
    "</w:document>" +
    "</pkg:xmlData>" +
    "</pkg:part>" +
    "</pkg:package>"

$r.InsertXML($xmlFragment) | Out-Null
